{"js": "// Office.js (Word JavaScript API) edit script.\n// Updates the 7 existing paragraphs' text (date, title, body paragraphs 1-4),\n// then inserts 3 new paragraphs before the final URL paragraph and updates\n// that URL paragraph's text too, matching the target diff.\n\nconst afterTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -09.11.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"CROSS-ENTROPY IS ALL YOU NEED TO INVERT THE DATA GENERATING PROCESS\",\n  \"\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e9\u05da \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e1\u05e7\u05e8\u05ea\u05d9 \u05d0\u05ea\u05de\u05d5\u05dc \u05e9\u05d4\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05dc-SSL \u05d0\u05d5 Self-Supervised Learning \u05d1\u05d0\u05e8\u05d9\u05db\u05d5\u05ea. \u05de\u05d8\u05e8\u05ea SSL \u05d4\u05d9\u05d0 \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05e4\u05d9\u05e7 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d3\u05d0\u05d8\u05d4 \u05e2\u05d5\u05e6\u05de\u05ea\u05d9 \u05e9\u05d9\u05d4\u05d9\u05d4 \u05e7\u05dc \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05de\u05e0\u05d5 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd downstream \u05dc\u05d1\u05d9\u05e6\u05d5\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e2\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d6\u05d4 \u05d1\u05ea\u05d5\u05e8 backbone (\u05dc\u05de\u05e9\u05dc \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d5\u05e1\u05e4\u05ea \u05e9\u05db\u05d1\u05d5\u05ea, LoRA, \u05d0\u05d3\u05e4\u05d8\u05e8\u05d9\u05dd \u05d0\u05d5 \u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05d0\u05d7\u05e8\u05d5\u05ea \u05d4\u05d1\u05e0\u05d5\u05d9\u05d5\u05ea \u05e2\u05dc \u05d4-backbone \u05d4\u05d6\u05d4). \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d6\u05d4 \u05e6\u05e8\u05d9\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05d6\u05e7\u05e7 \u05d0\u05ea \u05db\u05dc \u05d4\u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d4\u05de\u05d4\u05d5\u05ea\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d6\u05d4 \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d3\u05d7\u05d5\u05e1\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d9\u05e2\u05d9\u05dc\u05d4.\",\n  \"\u05de\u05e9\u05d9\u05de\u05ea downstream \u05d4\u05e4\u05e9\u05d5\u05d8\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05d9\u05d0 \u05de\u05e9\u05d9\u05de\u05ea \u05e1\u05d9\u05d5\u05d5\u05d2 \u05d5\u05d1\u05de\u05e7\u05e8\u05d4 \u05d4\u05d6\u05d4 \u05de\u05d5\u05d3\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d8\u05d5\u05d1 \u05e6\u05e8\u05d9\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05d1\u05d9\u05df \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d9\u05d9\u05da \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea (\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05de\u05d0\u05d5\u05de\u05df \u05e2\u05dc \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0 \u05de\u05ea\u05d5\u05d9\u05d2). \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05dc \u05d0\u05ea\u05de\u05d5\u05dc \u05d4\u05e6\u05d9\u05e2 \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e9\u05d9\u05d5\u05d3\u05e2 \u05dc\u05d6\u05d4\u05d5\u05ea \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc\u05d4. \u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05dc \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05de\u05e7\u05d1\u05dc\u05ea \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 \u05de\u05e9\u05dc\u05d4 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e2\u05dd 10L \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d0\u05d6 \u05d9\u05e9 \u05dc\u05e0\u05d5 10K \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea). \u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d1\u05e0\u05d5\u05e1\u05e3 \u05dc\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 (\u05de\u05d5\u05d3\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2) \u05e9\u05de\u05de\u05e4\u05d4 (\u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea) \u05d0\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e2\u05dd \u05dc\u05d5\u05e1 cross-entropy. \",\n  \"\u05d0\u05d6 \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05dc \u05d0\u05ea\u05de\u05d5\u05dc \u05d8\u05e2\u05df \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d2\u05d9\u05e2 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d7\u05d6\u05e7\u05d9\u05dd \u05e2\u05dd \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d5 (\u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea downstream \u05de\u05e1\u05d5\u05d2 \u05e1\u05d9\u05d5\u05d5\u05d2) \u05d5\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05d4\u05d5\u05db\u05d9\u05d7 \u05db\u05de\u05d4 \u05d8\u05e2\u05e0\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05e9\u05e0\u05d3\u05d5\u05df \u05d1\u05de\u05d0\u05de\u05e8 (\u05d8\u05d5\u05d1 \u05d6\u05d4 \u05dc\u05d0 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d0\u05d1\u05dc \u05e7\u05e8\u05d5\u05d1) \u05e9\u05e1\u05e7\u05e8\u05e0\u05d5 \u05d0\u05ea\u05de\u05d5\u05dc \u05ea\u05d7\u05ea \u05d4\u05e0\u05d7\u05d5\u05ea \u05d3\u05d9 \u05d4\u05d2\u05d9\u05d5\u05e0\u05d9\u05d5\u05ea. \u05d4\u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05ea\u05de\u05d8\u05d9 \u05d5\u05d0\u05e0\u05e1\u05d4 \u05dc\u05d4\u05e1\u05d1\u05d9\u05e8 \u05d0\u05ea \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9 \u05d1\u05dc\u05d9 \u05dc\u05e6\u05dc\u05d5\u05dc \u05dc\u05e0\u05d5\u05e1\u05d7\u05d0\u05d5\u05ea \u05d5\u05dc\u05dc\u05d0 \u05d4\u05ea\u05e2\u05de\u05e7\u05d5\u05d9\u05d5\u05ea \u05d9\u05ea\u05e8 \u05dc\u05e4\u05e8\u05d8\u05d9\u05dd \u05de\u05ea\u05de\u05d8\u05d9\u05d9\u05dd \u05dc\u05d0 \u05de\u05d4\u05d5\u05ea\u05d9\u05d9\u05dd.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e0\u05d9\u05d7\u05d9\u05dd \u05db\u05de\u05d4 \u05d4\u05e0\u05d7\u05d5\u05ea \u05e9\u05e2\u05d5\u05d6\u05e8\u05d5\u05ea \u05dc\u05d4\u05dd \u05dc\u05d7\u05e7\u05d5\u05e8 \u05d0\u05ea \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d6\u05d5. \u05d4\u05d4\u05e0\u05d7\u05d4 \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05de\u05e0\u05d9\u05d7\u05d4 \u05e9\u05d9\u05e9 \u05ea\u05d4\u05dc\u05d9\u05da \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9 \u05d4\u05de\u05d2\u05e0\u05e8\u05d8 \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d4\u05e9\u05d9\u05d9\u05db\u05d9\u05dd \u05dc\u05db\u05de\u05d4 \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea (\u05de\u05e1\u05e4\u05e8\u05dd \u05d9\u05d3\u05d5\u05e2). \u05d1\u05e4\u05e8\u05d8 \u05d4\u05d9\u05d0 \u05de\u05d3\u05d1\u05e8\u05ea \u05e2\u05dc \u05db\u05da \u05e9\u05e7\u05d9\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9 g \u05d4\u05de\u05d2\u05e0\u05e8\u05d8 \u05d3\u05d0\u05d8\u05d4 \u05de\u05d9\u05d9\u05e6\u05d5\u05d2\u05d5 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 z. \u05d4\u05de\u05e9\u05ea\u05e0\u05d4 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 z \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 C \u05de\u05d5\u05d2\u05e8\u05dc \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea von Neumann-Fisher \u05d0\u05d5 vMF \u05d1\u05e7\u05e6\u05e8\u05d4. vMF \u05d4\u05d9\u05d0 \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e8\u05d1 \u05de\u05d9\u05de\u05d3\u05d9\u05ea \u05e2\u05dc \u05e1\u05e4\u05d9\u05e8\u05d4 \u05d1\u05e2\u05dc\u05ea \u05e8\u05d3\u05d9\u05d5\u05e1 \u05d0\u05d7\u05ea \u05d4\u05de\u05d5\u05d2\u05d3\u05e8\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d5\u05e7\u05d8\u05d5\u05e8 z_c \u05ea\u05d5\u05d7\u05dc\u05ea \u05d5\u05e4\u05e8\u05de\u05d8\u05e8 \u05e8\u05d9\u05db\u05d5\u05d6 (\u05e1\u05e7\u05dc\u05e8 \u05d4\u05de\u05d2\u05d3\u05d9\u05e8 \u05d0\u05ea \u05de\u05d9\u05d3\u05ea \u05d4\u05de\u05e8\u05d9\u05d7\u05d5\u05ea \u05e9\u05dc \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea). \",\n  \"\u05e2\u05db\u05e9\u05d9\u05d5 \u05d4\u05de\u05e9\u05e4\u05d8 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d1\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05d0\u05dd \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2 f (\u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8) \u05e2\u05dc\u05d9 \u05d9\u05d3\u05d9 \u05de\u05e7\u05e1\u05d5\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05d3\u05d5\u05de\u05d4 \u05dc\u05d6\u05d0\u05ea \u05de\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e7\u05d5\u05d3\u05dd \u05e8\u05e7 \u05e9\u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d9\u05d4\u05d9\u05d4 \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4(\u05d4\u05de\u05d9\u05d5\u05e6\u05d2\u05d5\u05ea \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9) \u05d5\u05dc\u05d0 \u05db\u05dc \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d9\u05d9\u05db\u05ea \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 \u05de\u05e9\u05dc\u05d4(\u05e0\u05db\u05d5\u05df \u05d6\u05d4 \u05dc\u05d0 \u05d0\u05d5\u05ea\u05d5 \u05d4\u05d3\u05d1\u05e8 \u05d0\u05d1\u05dc \u05e2\u05d3\u05d9\u05d9\u05df), \u05d9\u05e9 \u05e4\u05d9\u05e8\u05d5\u05e9 \u05d3\u05d9 \u05d9\u05e4\u05d4 \u05dc\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd w \u05d4\u05de\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05ea W \u05e9\u05d4\u05d9\u05d0 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9 \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 \u05e9\u05d0\u05e0\u05d5 \u05dc\u05d5\u05de\u05d3\u05d9\u05dd \u05de\u05d4\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05dc\u05de\u05e8\u05d7\u05d1 \u05d4\u05d3\u05d0\u05d8\u05d4. \",\n  \"\u05d1\u05de\u05e7\u05e8\u05d4 \u05d4\u05e4\u05e9\u05d5\u05d8 - \u05de\u05e9\u05e4\u05d8 \u05d0\u05d7\u05d3 \u05de\u05d2\u05d3\u05d9\u05e8 4 \u05de\u05e7\u05e8\u05d9\u05dd, \u05d4\u05ea\u05dc\u05d5\u05d9\u05d9\u05dd \u05d4\u05d0\u05dd \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd w (\u05d4\u05de\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d0\u05ea W) \u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d0\u05d7\u05e8\u05d9 (f(x, \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 w \u05de\u05d4\u05d5\u05d5\u05d9\u05dd \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e6\u05d9\u05d4 \u05d0\u05d5\u05e8\u05ea\u05d5\u05d2\u05d5\u05e0\u05dc\u05d9\u05ea \u05e9\u05dc \u05de\u05e8\u05db\u05d6\u05d9 \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea z_c \u05e9\u05de\u05de\u05e0\u05d5 \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05de\u05d5\u05d2\u05e8\u05dc\u05d9\u05dd (\u05db\u05dc\u05d5\u05de\u05e8 \u05d6\u05d4 \u05d0\u05d5\u05ea\u05dd \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05ea\u05d7\u05ea \u05e1\u05d9\u05d1\u05d5\u05d1 \u05e8\u05d1 \u05de\u05d9\u05de\u05d3\u05d9 \u05db\u05dc\u05e9\u05d4\u05d5). \u05db\u05dc\u05d5\u05de\u05e8 \u05e7\u05d9\u05d1\u05dc\u05e0\u05d5 w_i \u05e2\u05dd \u05de\u05d0\u05d5\u05d3 \u05e7\u05e9\u05d5\u05e8\u05d9\u05dd \u05dc\u05de\u05d1\u05e0\u05d4 \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4. \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d1\u05de\u05e7\u05e8\u05d4 \u05d4\u05d6\u05d4 \u05d4\u05d4\u05e8\u05db\u05d1\u05d4 \u05e9\u05dc \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 f (\u05de\u05d4 \u05e9\u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd) \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 g \u05d4\u05d9\u05e0\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e6\u05dc\u05d7\u05e0\u05d5 \u05dc\u05de\u05e6\u05d5\u05d0 \u05d0\u05ea \u05d4\u05d4\u05d5\u05e4\u05db\u05d9\u05ea \u05e9\u05dc \u05d4\u05d2\u05e0\u05e8\u05d8\u05d5\u05e8 g - \u05d5\u05d6\u05d4 \u05ea\u05d5\u05e6\u05d0\u05d4 \u05d3\u05d9 \u05d7\u05d6\u05e7\u05d4 (\u05de\u05e9\u05e4\u05d8 2 \u05de\u05e0\u05e1\u05d7 \u05d0\u05ea \u05d6\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9 \u05d8\u05d5\u05d1\u05d4).\",\n  \"\u05d4\u05d4\u05d5\u05db\u05d7\u05d5\u05ea \u05dc\u05d0 \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d1\u05db\u05dc\u05dc \u05d5\u05e2\u05dd \u05d6\u05d0\u05ea \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d7\u05e9\u05d5\u05d1 \u05d5\u05d0\u05e0\u05d9 \u05de\u05e7\u05d5\u05d5\u05d4 \u05e9\u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05e4\u05d7\u05d5\u05ea \u05dc\u05d4\u05e1\u05d1\u05d9\u05e8 \u05dc\u05db\u05dd \u05d0\u05ea \u05de\u05d4\u05d5\u05ea\u05d5.\",\n  \"https://arxiv.org/abs/2410.21869\"\n];\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\nif (paragraphs.length !== 7) {\n  throw new Error(\"Unexpected paragraph count: \" + paragraphs.length);\n}\n\n// Update the first 6 paragraphs in place (indices 0-5 map directly to afterTexts 0-5).\nfor (let i = 0; i < 6; i++) {\n  paragraphs[i].insertText(afterTexts[i], Word.InsertLocation.replace);\n}\n\n// The 7th paragraph (index 6) currently holds the arXiv link. Insert the three\n// new paragraphs (afterTexts[6..8]) immediately before it, then replace its\n// own text with the new link (afterTexts[9]).\nconst linkParagraph = paragraphs[6];\nfor (let i = 6; i <= 8; i++) {\n  linkParagraph.insertParagraph(afterTexts[i], Word.InsertLocation.before);\n}\nlinkParagraph.insertText(afterTexts[9], Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Updates the text of the 7 existing paragraphs (date, title, body 1-4),\n# then inserts 3 new paragraphs before the final URL paragraph and\n# updates that paragraph's text to the new arXiv link, matching the diff.\n\n$d = $word.ActiveDocument\n\n$afterTexts = @(\n  '\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -09.11.24: \u26a1\ufe0f\ud83d\ude80',\n  'CROSS-ENTROPY IS ALL YOU NEED TO INVERT THE DATA GENERATING PROCESS',\n  '\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e9\u05da \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e1\u05e7\u05e8\u05ea\u05d9 \u05d0\u05ea\u05de\u05d5\u05dc \u05e9\u05d4\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05dc-SSL \u05d0\u05d5 Self-Supervised Learning \u05d1\u05d0\u05e8\u05d9\u05db\u05d5\u05ea. \u05de\u05d8\u05e8\u05ea SSL \u05d4\u05d9\u05d0 \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d5\u05d3\u05dc \u05d4\u05de\u05e4\u05d9\u05e7 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d3\u05d0\u05d8\u05d4 \u05e2\u05d5\u05e6\u05de\u05ea\u05d9 \u05e9\u05d9\u05d4\u05d9\u05d4 \u05e7\u05dc \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05de\u05e0\u05d5 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd downstream \u05dc\u05d1\u05d9\u05e6\u05d5\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e2\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d6\u05d4 \u05d1\u05ea\u05d5\u05e8 backbone (\u05dc\u05de\u05e9\u05dc \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d5\u05e1\u05e4\u05ea \u05e9\u05db\u05d1\u05d5\u05ea, LoRA, \u05d0\u05d3\u05e4\u05d8\u05e8\u05d9\u05dd \u05d0\u05d5 \u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05d0\u05d7\u05e8\u05d5\u05ea \u05d4\u05d1\u05e0\u05d5\u05d9\u05d5\u05ea \u05e2\u05dc \u05d4-backbone \u05d4\u05d6\u05d4). \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d6\u05d4 \u05e6\u05e8\u05d9\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05d6\u05e7\u05e7 \u05d0\u05ea \u05db\u05dc \u05d4\u05ea\u05db\u05d5\u05e0\u05d5\u05ea \u05d4\u05de\u05d4\u05d5\u05ea\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d4\u05d6\u05d4 \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d3\u05d7\u05d5\u05e1\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d9\u05e2\u05d9\u05dc\u05d4.',\n  '\u05de\u05e9\u05d9\u05de\u05ea downstream \u05d4\u05e4\u05e9\u05d5\u05d8\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05d9\u05d0 \u05de\u05e9\u05d9\u05de\u05ea \u05e1\u05d9\u05d5\u05d5\u05d2 \u05d5\u05d1\u05de\u05e7\u05e8\u05d4 \u05d4\u05d6\u05d4 \u05de\u05d5\u05d3\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d8\u05d5\u05d1 \u05e6\u05e8\u05d9\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05d1\u05d9\u05df \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d9\u05d9\u05da \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea (\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05de\u05d0\u05d5\u05de\u05df \u05e2\u05dc \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0 \u05de\u05ea\u05d5\u05d9\u05d2). \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05dc \u05d0\u05ea\u05de\u05d5\u05dc \u05d4\u05e6\u05d9\u05e2 \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e9\u05d9\u05d5\u05d3\u05e2 \u05dc\u05d6\u05d4\u05d5\u05ea \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc\u05d4. \u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05dc \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05de\u05e7\u05d1\u05dc\u05ea \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 \u05de\u05e9\u05dc\u05d4 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e2\u05dd 10L \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d0\u05d6 \u05d9\u05e9 \u05dc\u05e0\u05d5 10K \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea). \u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e9\u05db\u05d1\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d1\u05e0\u05d5\u05e1\u05e3 \u05dc\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 (\u05de\u05d5\u05d3\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2) \u05e9\u05de\u05de\u05e4\u05d4 (\u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea) \u05d0\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e2\u05dd \u05dc\u05d5\u05e1 cross-entropy. ',\n  '\u05d0\u05d6 \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05dc \u05d0\u05ea\u05de\u05d5\u05dc \u05d8\u05e2\u05df \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d2\u05d9\u05e2 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d7\u05d6\u05e7\u05d9\u05dd \u05e2\u05dd \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d5 (\u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea downstream \u05de\u05e1\u05d5\u05d2 \u05e1\u05d9\u05d5\u05d5\u05d2) \u05d5\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05d4\u05d5\u05db\u05d9\u05d7 \u05db\u05de\u05d4 \u05d8\u05e2\u05e0\u05d5\u05ea \u05dc\u05d2\u05d1\u05d9 \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05e9\u05e0\u05d3\u05d5\u05df \u05d1\u05de\u05d0\u05de\u05e8 (\u05d8\u05d5\u05d1 \u05d6\u05d4 \u05dc\u05d0 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d0\u05d1\u05dc \u05e7\u05e8\u05d5\u05d1) \u05e9\u05e1\u05e7\u05e8\u05e0\u05d5 \u05d0\u05ea\u05de\u05d5\u05dc \u05ea\u05d7\u05ea \u05d4\u05e0\u05d7\u05d5\u05ea \u05d3\u05d9 \u05d4\u05d2\u05d9\u05d5\u05e0\u05d9\u05d5\u05ea. \u05d4\u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05ea\u05de\u05d8\u05d9 \u05d5\u05d0\u05e0\u05e1\u05d4 \u05dc\u05d4\u05e1\u05d1\u05d9\u05e8 \u05d0\u05ea \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9 \u05d1\u05dc\u05d9 \u05dc\u05e6\u05dc\u05d5\u05dc \u05dc\u05e0\u05d5\u05e1\u05d7\u05d0\u05d5\u05ea \u05d5\u05dc\u05dc\u05d0 \u05d4\u05ea\u05e2\u05de\u05e7\u05d5\u05d9\u05d5\u05ea \u05d9\u05ea\u05e8 \u05dc\u05e4\u05e8\u05d8\u05d9\u05dd \u05de\u05ea\u05de\u05d8\u05d9\u05d9\u05dd \u05dc\u05d0 \u05de\u05d4\u05d5\u05ea\u05d9\u05d9\u05dd.',\n  '\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e0\u05d9\u05d7\u05d9\u05dd \u05db\u05de\u05d4 \u05d4\u05e0\u05d7\u05d5\u05ea \u05e9\u05e2\u05d5\u05d6\u05e8\u05d5\u05ea \u05dc\u05d4\u05dd \u05dc\u05d7\u05e7\u05d5\u05e8 \u05d0\u05ea \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d6\u05d5. \u05d4\u05d4\u05e0\u05d7\u05d4 \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05de\u05e0\u05d9\u05d7\u05d4 \u05e9\u05d9\u05e9 \u05ea\u05d4\u05dc\u05d9\u05da \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9 \u05d4\u05de\u05d2\u05e0\u05e8\u05d8 \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d4\u05e9\u05d9\u05d9\u05db\u05d9\u05dd \u05dc\u05db\u05de\u05d4 \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea (\u05de\u05e1\u05e4\u05e8\u05dd \u05d9\u05d3\u05d5\u05e2). \u05d1\u05e4\u05e8\u05d8 \u05d4\u05d9\u05d0 \u05de\u05d3\u05d1\u05e8\u05ea \u05e2\u05dc \u05db\u05da \u05e9\u05e7\u05d9\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9 g \u05d4\u05de\u05d2\u05e0\u05e8\u05d8 \u05d3\u05d0\u05d8\u05d4 \u05de\u05d9\u05d9\u05e6\u05d5\u05d2\u05d5 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 z. \u05d4\u05de\u05e9\u05ea\u05e0\u05d4 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 z \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 C \u05de\u05d5\u05d2\u05e8\u05dc \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea von Neumann-Fisher \u05d0\u05d5 vMF \u05d1\u05e7\u05e6\u05e8\u05d4. vMF \u05d4\u05d9\u05d0 \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e8\u05d1 \u05de\u05d9\u05de\u05d3\u05d9\u05ea \u05e2\u05dc \u05e1\u05e4\u05d9\u05e8\u05d4 \u05d1\u05e2\u05dc\u05ea \u05e8\u05d3\u05d9\u05d5\u05e1 \u05d0\u05d7\u05ea \u05d4\u05de\u05d5\u05d2\u05d3\u05e8\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d5\u05e7\u05d8\u05d5\u05e8 z_c \u05ea\u05d5\u05d7\u05dc\u05ea \u05d5\u05e4\u05e8\u05de\u05d8\u05e8 \u05e8\u05d9\u05db\u05d5\u05d6 (\u05e1\u05e7\u05dc\u05e8 \u05d4\u05de\u05d2\u05d3\u05d9\u05e8 \u05d0\u05ea \u05de\u05d9\u05d3\u05ea \u05d4\u05de\u05e8\u05d9\u05d7\u05d5\u05ea \u05e9\u05dc \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea). ',\n  '\u05e2\u05db\u05e9\u05d9\u05d5 \u05d4\u05de\u05e9\u05e4\u05d8 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d1\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05d0\u05dd \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d9\u05d9\u05e6\u05d5\u05d2 f (\u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8) \u05e2\u05dc\u05d9 \u05d9\u05d3\u05d9 \u05de\u05e7\u05e1\u05d5\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05d3\u05d5\u05de\u05d4 \u05dc\u05d6\u05d0\u05ea \u05de\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e7\u05d5\u05d3\u05dd \u05e8\u05e7 \u05e9\u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05d9\u05d4\u05d9\u05d4 \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4(\u05d4\u05de\u05d9\u05d5\u05e6\u05d2\u05d5\u05ea \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9) \u05d5\u05dc\u05d0 \u05db\u05dc \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d9\u05d9\u05db\u05ea \u05dc\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4 \u05de\u05e9\u05dc\u05d4(\u05e0\u05db\u05d5\u05df \u05d6\u05d4 \u05dc\u05d0 \u05d0\u05d5\u05ea\u05d5 \u05d4\u05d3\u05d1\u05e8 \u05d0\u05d1\u05dc \u05e2\u05d3\u05d9\u05d9\u05df), \u05d9\u05e9 \u05e4\u05d9\u05e8\u05d5\u05e9 \u05d3\u05d9 \u05d9\u05e4\u05d4 \u05dc\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd w \u05d4\u05de\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05ea W \u05e9\u05d4\u05d9\u05d0 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9 \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 \u05e9\u05d0\u05e0\u05d5 \u05dc\u05d5\u05de\u05d3\u05d9\u05dd \u05de\u05d4\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05dc\u05de\u05e8\u05d7\u05d1 \u05d4\u05d3\u05d0\u05d8\u05d4. ',\n  '\u05d1\u05de\u05e7\u05e8\u05d4 \u05d4\u05e4\u05e9\u05d5\u05d8 - \u05de\u05e9\u05e4\u05d8 \u05d0\u05d7\u05d3 \u05de\u05d2\u05d3\u05d9\u05e8 4 \u05de\u05e7\u05e8\u05d9\u05dd, \u05d4\u05ea\u05dc\u05d5\u05d9\u05d9\u05dd \u05d4\u05d0\u05dd \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd w (\u05d4\u05de\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d0\u05ea W) \u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d0\u05d7\u05e8\u05d9 (f(x, \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 w \u05de\u05d4\u05d5\u05d5\u05d9\u05dd \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e6\u05d9\u05d4 \u05d0\u05d5\u05e8\u05ea\u05d5\u05d2\u05d5\u05e0\u05dc\u05d9\u05ea \u05e9\u05dc \u05de\u05e8\u05db\u05d6\u05d9 \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea z_c \u05e9\u05de\u05de\u05e0\u05d5 \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05de\u05d5\u05d2\u05e8\u05dc\u05d9\u05dd (\u05db\u05dc\u05d5\u05de\u05e8 \u05d6\u05d4 \u05d0\u05d5\u05ea\u05dd \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05ea\u05d7\u05ea \u05e1\u05d9\u05d1\u05d5\u05d1 \u05e8\u05d1 \u05de\u05d9\u05de\u05d3\u05d9 \u05db\u05dc\u05e9\u05d4\u05d5). \u05db\u05dc\u05d5\u05de\u05e8 \u05e7\u05d9\u05d1\u05dc\u05e0\u05d5 w_i \u05e2\u05dd \u05de\u05d0\u05d5\u05d3 \u05e7\u05e9\u05d5\u05e8\u05d9\u05dd \u05dc\u05de\u05d1\u05e0\u05d4 \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4. \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d1\u05de\u05e7\u05e8\u05d4 \u05d4\u05d6\u05d4 \u05d4\u05d4\u05e8\u05db\u05d1\u05d4 \u05e9\u05dc \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 f (\u05de\u05d4 \u05e9\u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd) \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 g \u05d4\u05d9\u05e0\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e6\u05dc\u05d7\u05e0\u05d5 \u05dc\u05de\u05e6\u05d5\u05d0 \u05d0\u05ea \u05d4\u05d4\u05d5\u05e4\u05db\u05d9\u05ea \u05e9\u05dc \u05d4\u05d2\u05e0\u05e8\u05d8\u05d5\u05e8 g - \u05d5\u05d6\u05d4 \u05ea\u05d5\u05e6\u05d0\u05d4 \u05d3\u05d9 \u05d7\u05d6\u05e7\u05d4 (\u05de\u05e9\u05e4\u05d8 2 \u05de\u05e0\u05e1\u05d7 \u05d0\u05ea \u05d6\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9 \u05d8\u05d5\u05d1\u05d4).',\n  '\u05d4\u05d4\u05d5\u05db\u05d7\u05d5\u05ea \u05dc\u05d0 \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d1\u05db\u05dc\u05dc \u05d5\u05e2\u05dd \u05d6\u05d0\u05ea \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d7\u05e9\u05d5\u05d1 \u05d5\u05d0\u05e0\u05d9 \u05de\u05e7\u05d5\u05d5\u05d4 \u05e9\u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05e4\u05d7\u05d5\u05ea \u05dc\u05d4\u05e1\u05d1\u05d9\u05e8 \u05dc\u05db\u05dd \u05d0\u05ea \u05de\u05d4\u05d5\u05ea\u05d5.',\n  'https://arxiv.org/abs/2410.21869'\n)\n\nif ($d.Paragraphs.Count -ne 7) {\n  throw \"Unexpected paragraph count: $($d.Paragraphs.Count)\"\n}\n\n# Update the first 6 paragraphs in place (indices 1-6 map to afterTexts[0..5]).\nfor ($i = 1; $i -le 6; $i++) {\n  $d.Paragraphs($i).Range.Text = $afterTexts[$i - 1]\n}\n\n# Paragraph 7 currently holds the arXiv link. Insert the three new\n# paragraphs (afterTexts[6..8]) immediately before it, then replace its\n# own text with the new link (afterTexts[9]).\n#\n# InsertParagraphBefore() splits a new empty paragraph in just before the\n# paragraph at $linkIndex; that new empty paragraph takes over $linkIndex\n# and the old paragraph (the link, or the previously inserted one) shifts\n# to $linkIndex + 1. So re-fetch by index each time and bump the index\n# after each insert to keep appending in the right reading order.\n$linkIndex = 7\nfor ($i = 6; $i -le 8; $i++) {\n  $d.Paragraphs($linkIndex).Range.InsertParagraphBefore()\n  $d.Paragraphs($linkIndex).Range.Text = $afterTexts[$i]\n  $linkIndex = $linkIndex + 1\n}\n$d.Paragraphs($linkIndex).Range.Text = $afterTexts[9]\n\nWrite-Output \"done\"\n"}
